$wb = $excel.ActiveWorkbook

# Rename sheets to mirror power sector naming system
$wsProd = $wb.Worksheets.Item("IFTQfS-PTC")
$wsProd.Name = "IFTQfS-production"

$wsInv = $wb.Worksheets.Item("IFTQfS-ITC")
$wsInv.Name = "IFTQfS-investment"

# Update the "About" sheet text to reflect new naming
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Range("A1").Value = "IFTQfS Industrial Fuel Types Qualifying for Clean Heat Production Subsidy"
$wsAbout.Range("A2").Value = "IFTQfS Industrial Fuel Types Qualifying for Clean Heat Investment Subsidy"
$wsAbout.Range("A8").Value = "This variable specifies which fuel types qualify for the clean industrial heat production,"
$wsAbout.Range("A9").Value = "clean industrial heat investment policies (when used for industrial heating purposes),"

$wsAbout.Range("A10").Select()

$wsInv.Activate()
$wsInv.Range("B14").Select()
